$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.676.75'
$ws.Range('E2').Value = '  -2.77%  '
$ws.Range('D3').Value = '1.985.53'
$ws.Range('E3').Value = '  -3.63%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''245.59'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D7').Value = '''59.60'
$ws.Range('E7').Value = '  +7.93%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '''59.13'
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('E11').Value = '  -1.68%  '
$ws.Range('E12').Value = '  -2.39%  '
$ws.Range('D13').Value = '''0.953'
$ws.Range('E13').Value = '  +1.06%  '
$ws.Range('D14').Value = '''14.67'
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('D15').Value = '2.271.82'
$ws.Range('E15').Value = '  -3.78%  '
$ws.Range('E16').Value = '  -2.72%  '
$ws.Range('D17').Value = '1.993.33'
$ws.Range('E17').Value = '  -3.16%  '
$ws.Range('D18').Value = '''18.73'
$ws.Range('E18').Value = '  +9.03%  '
$ws.Range('D19').Value = '35.592.26'
$ws.Range('E19').Value = '  -2.84%  '
$ws.Range('D20').Value = '''71.71'
$ws.Range('E20').Value = '  -0.64%  '
$ws.Range('D21').Value = '0.0₃0851'
$ws.Range('E21').Value = '  -1.71%  '
$ws.Range('D22').Value = '''5.23'
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').Value = '''233.38'
$ws.Range('E23').Value = '  -2.30%  '
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  +17.22%  '
$ws.Range('D26').Value = '''2.28'
$ws.Range('E26').Value = '  -4.41%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '''9.23'
$ws.Range('E27').Value = '  -1.29%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '''165.34'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  -4.62%  '
$ws.Range('E30').Value = '  -2.70%  '
$ws.Range('E31').Value = '  -4.43%  '
$ws.Range('E32').Value = '  -6.52%  '
$ws.Range('D33').Value = '''0.0958'
$ws.Range('E33').Value = '  +13.29%  '
$ws.Range('D34').Value = '''0.0599'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('D35').Value = '''2.46'
$ws.Range('E35').Value = '  +10.68%  '
$ws.Range('D36').Value = '''4.37'
$ws.Range('E36').Value = '  -3.27%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  -2.65%  '
$ws.Range('D39').Value = '''5.55'
$ws.Range('E39').Value = '  +10.63%  '
$ws.Range('E40').Value = '  -1.75%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.0214'
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').Value = '''2.85'
$ws.Range('E42').Value = '  -2.12%  '
$ws.Range('D43').Value = '''0.0925'
$ws.Range('E43').Value = '  +2.49%  '
$ws.Range('D44').Value = '''7.83'
$ws.Range('E44').Value = '  +1.48%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').Value = '''1.09'
$ws.Range('E45').Value = '  -1.54%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '''16.47'
$ws.Range('E46').Value = '  +1.98%  '
$ws.Range('D47').Value = '''93.78'
$ws.Range('E47').Value = '  -1.24%  '
$ws.Range('D48').Value = '1.366.20'
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('D49').Value = '''2.90'
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('D50').Value = '''47.04'
$ws.Range('E50').Value = '  +3.07%  '
$ws.Range('E51').Value = '  -0.16%  '
